$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vocabulary entries to append: word, definition, example1, example2, row height
$entries = @(
    @("folk", "a specific group of people, distinguished by common nationality, background, or style of life.", "most of rich folk lives here.", "Thanks to the folks at NBC.", 75),
    @("grumpy", "bad-tempered and easily annoyed", "you are grumpy old man.", "Come back and see me when you're less grumpy.", 45),
    @("grapple", "to fight or struggle with someone, holding them tightly.", "the goverment is grappling with inflation.", "Two men grappled with a guard at the door.", 45),
    @("inflation", "a continuing increase in prices, or the rate at which prices increase", "Inflation is now at over 16%.", "Too much government borrowing can lead to inflation.", 45),
    @("grudge", "a feeling of dislike for someone because you cannot forget that they harmed you in the past", "I always feel she holds a grudge against me", "Is there anyone who might have had a grudge against her?", 60)
)

$startRow = 84

for ($i = 0; $i -lt $entries.Length; $i++) {
    $row = $startRow + $i
    $entry = $entries[$i]

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]

    # Column A: wrap text, vertically centered (matches style used by existing word column)
    $colA = $ws.Range("A$row")
    $colA.Font.Name = "Arial"
    $colA.Font.Size = 12
    $colA.Font.Color = 8421504
    $colA.WrapText = $true
    $colA.VerticalAlignment = -4108

    # Columns B:D: wrap text, vertically top-aligned (matches style used by definition/example columns)
    $colBD = $ws.Range("B$row`:D$row")
    $colBD.Font.Name = "Arial"
    $colBD.Font.Size = 12
    $colBD.Font.Color = 8421504
    $colBD.WrapText = $true
    $colBD.VerticalAlignment = -4160

    $ws.Rows.Item($row).RowHeight = $entry[4]
}

$ws.Range("D94").Select()
